$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")
$ws.Activate()

# Halve the "Price Per Email" discount factors on row 14 (D:F), and set
# new discounted rates for the Platinum/Diamond/Ultimate tiers (H:J).
$ws.Range("D14").Value = 0.5
$ws.Range("E14").Value = 0.4
$ws.Range("F14").Value = 0.3
$ws.Range("H14").Value = 0.25
$ws.Range("I14").Value = 0.2
$ws.Range("J14").Value = 0.17

# New "Enemy" comparison row (22): ratio of competitor margin to 70%.
$ws.Range("C22").Value = "Enemy"
$ws.Range("D22").Formula = "=D27/70%"
$ws.Range("E22").Formula = "=E27/70%"
$ws.Range("F22").Formula = "=F27/70%"
$ws.Range("G22").Formula = "=G27/70%"
$ws.Range("H22").Formula = "=H27/70%"
$ws.Range("I22").Formula = "=I27/70%"
$ws.Range("J22").Formula = "=J27/70%"
$ws.Range("D22:J22").NumberFormat = "0.00"

# New "Datacellectief" label row (23).
$ws.Range("C23").Value = "Datacellectief"

# New competitor pricing data (rows 25:27).
$ws.Range("D25").Value = 500
$ws.Range("E25").Value = 1000
$ws.Range("F25").Value = 2500
$ws.Range("H25").Value = 5000
$ws.Range("I25").Value = 15000
$ws.Range("J25").Value = 20000

$ws.Range("D26").Value = 175
$ws.Range("E26").Value = 275
$ws.Range("F26").Value = 475
$ws.Range("H26").Value = 825
$ws.Range("I26").Value = 1800
$ws.Range("J26").Value = 2200

$ws.Range("D27").Formula = "=D26/D25"
$ws.Range("E27").Formula = "=E26/E25"
$ws.Range("F27").Formula = "=F26/F25"
$ws.Range("H27").Formula = "=H26/H25"
$ws.Range("I27").Formula = "=I26/I25"
$ws.Range("J27").Formula = "=J26/J25"

# Update selection / view to match the edited state.
$ws.Range("G18").Select() | Out-Null
